$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.87"
$ws.Range("E2").Value = "'-0.38%"
$ws.Range("D3").Value = "'43.56"
$ws.Range("E3").Value = "'4.21%"
$ws.Range("D4").Value = "'5.612"
$ws.Range("E4").Value = "'-1.33%"
$ws.Range("D5").Value = "'0.08201"
$ws.Range("E5").Value = "'-1.95%"
$ws.Range("D6").Value = "'8.764"
$ws.Range("E6").Value = "'-0.36%"
$ws.Range("E7").Value = "'-3.09%"
$ws.Range("E8").Value = "'-5.64%"
$ws.Range("D9").Value = "'2.846"
$ws.Range("E9").Value = "'-4.26%"
$ws.Range("D10").Value = "'0.9440"
$ws.Range("E10").Value = "'1.87%"
$ws.Range("D11").Value = "'0.1202"
$ws.Range("E11").Value = "'-7.02%"
$ws.Range("D12").Value = "'0.1918"
$ws.Range("E12").Value = "'-2.82%"
$ws.Range("D13").Value = "'0.09827"
$ws.Range("E13").Value = "'2.69%"
$ws.Range("D14").Value = "'0.04352"
$ws.Range("E14").Value = "'11.24%"
$ws.Range("E15").Value = "'0.89%"
$ws.Range("D16").Value = "'0.001279"
$ws.Range("E16").Value = "'-2.59%"
$ws.Range("D17").Value = "'0.006019"
$ws.Range("E17").Value = "'-1.12%"
$ws.Range("E18").Value = "'1.81%"
$ws.Range("D19").Value = "'0.3536"
$ws.Range("D20").Value = "'8.743"
$ws.Range("E20").Value = "'6.16%"
$ws.Range("D21").Value = "'0.1369"
$ws.Range("E21").Value = "'-0.18%"
$ws.Range("D22").Value = "'0.2523"
$ws.Range("E22").Value = "'4.63%"
$ws.Range("D23").Value = "'0.04389"
$ws.Range("E23").Value = "'-0.59%"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'-0.72%"
$ws.Range("D25").Value = "'0.004314"
$ws.Range("E25").Value = "'-1.40%"
$ws.Range("D26").Value = "'0.0001236"
$ws.Range("E26").Value = "'2.91%"
$ws.Range("D27").Value = "'0.0004009"
$ws.Range("E27").Value = "'31.64%"
$ws.Range("D39").Value = "'0.02781"
$ws.Range("E39").Value = "'-1.11%"
$ws.Range("D40").Value = "'0.05731"
$ws.Range("E40").Value = "'3.34%"
$ws.Range("D41").Value = "'0.007938"
$ws.Range("E41").Value = "'1.78%"
$ws.Range("D42").Value = "'0.009768"
$ws.Range("E42").Value = "'7.08%"
$ws.Range("D43").Value = "'0.1420"
$ws.Range("E43").Value = "'-0.98%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-1.74%"
$ws.Range("D45").Value = "'0.009660"
$ws.Range("E45").Value = "'-12.99%"
$ws.Range("D46").Value = "'0.00007360"
$ws.Range("E46").Value = "'4.47%"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.49%"
$ws.Range("D48").Value = "'0.003459"
$ws.Range("E48").Value = "'-1.61%"
$ws.Range("D49").Value = "'0.002281"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.49%"
$ws.Range("D51").Value = "'0.0002010"
$ws.Range("E51").Value = "'0.49%"
